# Daily attendance processing - 2026-01-08 19:09:52
#
# For each row in the "Recorded By" column (G), the value is a
# comma-separated list of users/systems that touched the record. This
# pass normalizes the ordering of the last two entries in that list:
# whenever the second-to-last entry sorts strictly before the last entry
# (case-sensitive / ordinal comparison), the two are swapped.

function OrdinalLess($ordA, $ordB) {
    $ordLa = $ordA.Length
    $ordLb = $ordB.Length
    $ordMinLen = [Math]::Min($ordLa, $ordLb)
    $ordI = 0
    while ($ordI -lt $ordMinLen) {
        $ordCa = [int][char]$ordA[$ordI]
        $ordCb = [int][char]$ordB[$ordI]
        if ($ordCa -lt $ordCb) { return $true }
        if ($ordCa -gt $ordCb) { return $false }
        $ordI = $ordI + 1
    }
    return $ordLa -lt $ordLb
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $txt = $cell.Value2
    if ($txt -eq $null) { continue }
    if ($txt -eq "") { continue }

    $parts = @($txt -split ", ")
    $n = $parts.Count
    if ($n -ge 2) {
        $last = $parts[$n - 1]
        $secondLast = $parts[$n - 2]
        if (OrdinalLess $secondLast $last) {
            $parts[$n - 2] = $last
            $parts[$n - 1] = $secondLast
            $newTxt = [string]::Join(", ", $parts)
            $cell.Value2 = $newTxt
        }
    }
}
